$wb = $excel.ActiveWorkbook

# --- 1) Rename header cells on the existing sheets -------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet at the end --------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the page-margin conventions used by the rest of the workbook
# (PageSetup margins are expressed in points; 1 inch = 72 points).
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# --- 3) Headers ---------------------------------------------------------
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- 4) Copy the date-formatted style down column A -------------------------
$wsWeekly.Range("A2:A7").Copy()
$wsForecast.Range("A2:A15").PasteSpecial(-4122)

# --- 5) Data rows ---------------------------------------------------------
$rows = @(
    @(44934.99999999999, 63, 20.68858672638984, 104.523448059829),
    @(44941.99999999999, 54, 15.4033287809745, 94.82555146770218),
    @(44948.99999999999, 44, 4.030025228323337, 84.68389786866298),
    @(44955.99999999999, 35, -4.517242755498605, 74.94871443549107),
    @(44962.99999999999, 25, -16.65744344887061, 65.15224293887653),
    @(44969.99999999999, 16, -26.00281904795503, 56.41763574223919),
    @(44976.99999999999, 6, -33.37356986738283, 48.70964637291363),
    @(44983.99999999999, 0, -42.75779835955561, 35.5306212436395),
    @(44990.99999999999, 0, -54.59376357176072, 28.95327251574026),
    @(44997.99999999999, 0, -62.033412715785, 17.75242676082814),
    @(45004.99999999999, 0, -69.86781898963893, 7.059598702824569),
    @(45011.99999999999, 0, -83.64812517140582, -0.9197426623889499),
    @(45018.99999999999, 0, -91.34582927729613, -12.3091231828835),
    @(45025.99999999999, 0, -100.9527481516957, -19.03600852751941)
)

$r = 2
foreach ($row in $rows) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
